$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectPlan")
$ws.Activate()

# Fill the "THÀNH VIÊN" (member) column (E) for rows 12-16 with "Long"
$ws.Range("E12:E16").Value = "Long"

# Fill the "THÀNH VIÊN" (member) column (E) for rows 17-20 with "Hien"
$ws.Range("E17:E20").Value = "Hien"

# Update the view: scroll so row 7 is at top of the window, and select G17
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G17").Select()
